# Commit: "change advances to have A"
# In the "codes" sheet, the rows describing "Advance" entries had their
# column D (the "Code" used when an on-field advance is recorded) updated
# so several of the short codes are now prefixed with "A" (for Advance).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("codes")

$ws.Range("D41").Value = "AFC"   # Fielder's Choice:         FC  -> AFC
$ws.Range("D37").Value = "AE"    # Error:                    E   -> AE
$ws.Range("D38").Value = "AE"    # Error after Passed Ball:  E   -> AE
$ws.Range("D39").Value = "AE"    # Error after Pickoff:      E   -> AE
$ws.Range("D40").Value = "AE"    # Error after Steal:        E   -> AE
$ws.Range("D43").Value = "APB"   # Passed Ball:               PB  -> APB
$ws.Range("D46").Value = "AWP"   # Wild Pitch:                WP  -> AWP
$ws.Range("D34").Value = "AX"    # Unknown:                 X   -> AX
$ws.Range("D47").Value = "AE"    # Error after Wild Pitch:   E   -> AE

# Restore the selection state that was active when the workbook was saved.
$ws.Range("I13").Select()
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 10
$activeWindow.ScrollColumn = 4
